$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.073.85'
$ws.Range("E2").Value = '  -1.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.792.55'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.46'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.38'
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("E9").Value = '  -1.89%  '
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0931'
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.048.32'
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.793.18'
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.84'
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("E15").Value = '  -2.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.050.76'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("E17").Value = '  -3.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.94'
$ws.Range("E18").Value = '  -1.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.12'
$ws.Range("E19").Value = '  -3.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0783'
$ws.Range("E20").Value = '  -2.50%  '
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.67'
$ws.Range("E22").Value = '  -3.44%  '
$ws.Range("E23").Value = '  -4.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.00'
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("E26").Value = '  -0.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.00'
$ws.Range("E27").Value = '  -2.10%  '
$ws.Range("E28").Value = '  -2.09%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0517'
$ws.Range("E30").Value = '  -1.97%  '
$ws.Range("E31").Value = '  +0.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.65'
$ws.Range("E32").Value = '  -3.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.49'
$ws.Range("E33").Value = '  -3.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.81'
$ws.Range("E34").Value = '  -4.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.384.62'
$ws.Range("E35").Value = '  -3.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.646'
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("E37").Value = '  -2.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0185'
$ws.Range("E38").Value = '  -3.64%  '
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.07'
$ws.Range("E41").Value = '  -3.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.913'
$ws.Range("E42").Value = '  -5.02%  '
$ws.Range("E44").Value = '  +8.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0496'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("E46").Value = '  -0.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '107.11'
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.949.55'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.83'
$ws.Range("E49").Value = '  -3.54%  '
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.88'
$ws.Range("E51").Value = '  -3.50%  '
